$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Is Active") currently stores text "True"/"False" as shared strings.
# Push up actual boolean values instead.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
